$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 50: date-like text in column A must stay text (not become a real
# date), so enter it with a leading apostrophe like a user forcing text.
$ws.Range("A50").Value = "'2025/10/02"
$ws.Range("B50").Value = "木"
$ws.Range("C50").Value = 13
$ws.Range("D50").Value = 3
